# Update gh-pages output: add "张家港·授渔4.5万圣夜观影节" event (new row 6) to the
# "展览" and "全部类型" sheets, bump a couple of "想去人数" counters, and fix the
# "想去人数" value recorded for a few already-listed events.

$wb = $excel.ActiveWorkbook

function Set-TextCell($ws, $row, $col, $text) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.ClearFormats()
}

function Insert-NewExpoRow($ws) {
    # Shift rows 6.. down by one and seed row 6 with the new event.
    $ws.Rows.Item(6).Insert()

    $src = $ws.Cells.Item(7, 1)
    $dst = $ws.Cells.Item(6, 1)
    $dst.Value = 5
    $dst.Font.Bold = $src.Font.Bold
    $dst.HorizontalAlignment = $src.HorizontalAlignment
    $dst.VerticalAlignment = $src.VerticalAlignment
    $dst.Borders.LineStyle = $src.Borders.LineStyle

    Set-TextCell $ws 6 2 "2024-10-26"
    Set-TextCell $ws 6 3 "张家港·授渔4.5万圣夜观影节"
    Set-TextCell $ws 6 4 "大新镇乐橙广场3楼 大新中影时代国际影城(乐橙广场店)"
    Set-TextCell $ws 6 5 "2024.10.26 14:00-10.26 21:00"
    $ws.Cells.Item(6, 6).Value = 1
    $ws.Cells.Item(6, 7).Value = 40
    Set-TextCell $ws 6 8 "https://show.bilibili.com/platform/detail.html?id=93516"
    Set-TextCell $ws 6 9 "//i2.hdslb.com/bfs/openplatform/202410/iNUOHGZs1729059598080.png"
}

# ---------------------------------------------------------------------------
# Sheet "展览"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Cells.Item(2, 6).Value = 814          # 恋与深空only 想去人数 813 -> 814

Insert-NewExpoRow $ws1

# Rows 7.. now hold the old rows 6.. shifted down by one; a handful of them
# also received a small "想去人数" correction alongside the shift.
$ws1.Cells.Item(8, 6).Value = 249           # 漫语堂动漫嘉年华 245 -> 249
$ws1.Cells.Item(10, 6).Value = 1022         # 创世次元兽装同人only展 1021 -> 1022
$ws1.Cells.Item(11, 6).Value = 14           # NK漫展全天趴 13 -> 14
$ws1.Cells.Item(12, 6).Value = 519          # 女神异闻录only同人展 518 -> 519
$ws1.Cells.Item(15, 6).Value = 12878        # COME IN JOY 动漫品牌国潮文化节 12871 -> 12878
$ws1.Cells.Item(17, 6).Value = 5290         # Good jump ACG元旦跨年盛典国潮文化节 5283 -> 5290

# ---------------------------------------------------------------------------
# Sheet "演出"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Cells.Item(2, 6).Value = 122            # 足太Penta生日会2024 117 -> 122

# ---------------------------------------------------------------------------
# Sheet "全部类型"
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Cells.Item(2, 6).Value = 814           # 恋与深空only 想去人数 813 -> 814

Insert-NewExpoRow $ws4

$ws4.Cells.Item(8, 6).Value = 249            # 漫语堂动漫嘉年华 245 -> 249
$ws4.Cells.Item(10, 6).Value = 1022          # 创世次元兽装同人only展 1021 -> 1022
$ws4.Cells.Item(11, 6).Value = 14            # NK漫展全天趴 13 -> 14
$ws4.Cells.Item(12, 6).Value = 519           # 女神异闻录only同人展 518 -> 519
$ws4.Cells.Item(15, 6).Value = 12878         # COME IN JOY 动漫品牌国潮文化节 12871 -> 12878
$ws4.Cells.Item(16, 6).Value = 122           # 足太Penta生日会2024 117 -> 122
$ws4.Cells.Item(19, 6).Value = 5290          # Good jump ACG元旦跨年盛典国潮文化节 5283 -> 5290
